# Weekly update: insert a new price record for "Comercializadora del Agro de
# Limarí - Haba" above the existing row 56, shifting all subsequent rows down
# by one (dimension grows from A1:R75 to A1:R76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 56; everything currently at row 56
# onward (through 75) moves down to 57..76.
$ws.Rows("56:56").Insert()

# Populate the new row 56 with the latest weekly price record.
$ws.Range("A56").Value = 2
$ws.Range("B56").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44846
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100112026
$ws.Range("G56").Value = "Haba"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 1600
$ws.Range("K56").Value = 4500
$ws.Range("L56").Value = 5000
$ws.Range("M56").Value = 4750
$ws.Range("N56").Value = "$/saco 25 kilos"
$ws.Range("O56").Value = "Provincia de Limarí"
$ws.Range("P56").Value = 190
$ws.Range("Q56").Value = 25
$ws.Range("R56").Value = "Hortaliza"
